$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, L, M, N, O, P, Q, R, S, T)
# computed from the target diff (rows effectively rotate their data among
# each other - this reproduces the exact resulting cell values).

$rows = @{
    2  = @{ D = 44742; L = "Segunda"; M = 100; N = 14000; O = 15000; P = 14500; Q = "$/caja 18 kilos granel";    R = "Región de O'Higgins"; S = 806;   T = 18 }
    5  = @{ D = 44707; L = "Primera"; M = 60;  N = 12000; O = 13000; P = 12500; Q = "$/caja 12 kilos empedrada"; R = "Provincia de Curicó";  S = 1042;  T = 12 }
    6  = @{ D = 44334; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "$/caja 12 kilos granel";    R = "Región de O'Higgins"; S = 11500; T = 1 }
    7  = @{ D = 45084; L = "Primera"; M = 100; N = 17000; O = 18000; P = 17500; Q = "$/caja 18 kilos granel";    R = "Región del Maule";      S = 972;   T = 18 }
    8  = @{ D = 44708; L = "Primera"; M = 70;  N = 12000; O = 13000; P = 12571; Q = "$/caja 12 kilos empedrada"; R = "Provincia de Curicó";  S = 1048;  T = 12 }
    9  = @{ D = 44330; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "$/caja 18 kilos granel";    R = "Provincia de Curicó";  S = 861;   T = 18 }
    10 = @{ D = 44719; L = "Primera"; M = 50;  N = 14000; O = 15000; P = 14400; Q = "$/caja 18 kilos granel";    R = "Región del Maule";      S = 800;   T = 18 }
    11 = @{ D = 44714; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "$/caja 18 kilos granel";    R = "Región de O'Higgins"; S = 806;   T = 18 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Range("D$r").Value = $v.D
    $ws.Range("L$r").Value = $v.L
    $ws.Range("M$r").Value = $v.M
    $ws.Range("N$r").Value = $v.N
    $ws.Range("O$r").Value = $v.O
    $ws.Range("P$r").Value = $v.P
    $ws.Range("Q$r").Value = $v.Q
    $ws.Range("R$r").Value = $v.R
    $ws.Range("S$r").Value = $v.S
    $ws.Range("T$r").Value = $v.T
}
